$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column K
$ws.Range("K1").Value = "intervention_type"

# Copy the style from the existing header cell J1 so K1 matches the other headers
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122) # xlPasteFormats

# Fill in intervention_type values for each clinical trial row
$values = @{
    2  = "BIOLOGICAL"
    3  = "OTHER"
    4  = "PROCEDURE"
    5  = "PROCEDURE"
    6  = "OTHER"
    7  = "OTHER"
    8  = "PROCEDURE"
    9  = "PROCEDURE"
    10 = "OTHER"
    11 = "BIOLOGICAL"
    12 = "PROCEDURE"
    13 = "BIOLOGICAL"
    14 = "PROCEDURE"
    15 = "DRUG"
    16 = "OTHER"
    17 = "DRUG"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}
